{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the text changes described by the diff:\n//  - \"Word Wide Web General\" bullets 1-3 reworded\n//  - \"Web Server\" bullet 1 extended, bullet 4 reworded\n//  - \"Web Client\" bullet 3 extended\n//  - \"Hello World\" bullets 4 & 5 reworded\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of exact current paragraph text -> new paragraph text.\n// Using the full current text as the search key keeps this robust even\n// though the source runs are split/merged in various ways.\nconst replacements = [\n  [\n    \"1.HTTP,CSS und Javascript haben WWW zum Erfolg gemacht.\",\n    \"1.HTTP, HTML und URL haben WWW zum Erfolg gemacht.\"\n  ],\n  [\n    \"2.IP-Adressen,URI oder DNS adressieren Ressourcen im WEB .\",\n    \"2.URL adressieren Ressourcen im WEB .\"\n  ],\n  [\n    \"3.Bestandteile einer URL: Verzeichnis,Dateiname,Textmarke,Username,Portnummer,CGI-Nummer,Protkoll und Servername.\",\n    \"3.Bestandteile einer URL: schema, special part, username, host, port, part, query, fragment\"\n  ],\n  [\n    \"1.Als Webserver bezeichnet man jene Server, die zur Verbreitung von Webinhalten im Internet dienen.\",\n    \"1.Als Webserver bezeichnet man jene Server, die Daten speichern und die zur Verbreitung von Webinhalten im Internet dienen.\"\n  ],\n  [\n    \"4.Durch ASP,PHP werden Inhalte in Webseiten beschrieben.\",\n    \"4.Durch HTML, PHP, ASHP werden Inhalte in Webseiten beschrieben.\"\n  ],\n  [\n    \"3.Web Client und Webbrowser sind das gleiche.\",\n    \"3.Web Client und Webbrowser sind nicht das gleiche.\"\n  ],\n  [\n    \"4 Mit 2x strong kann man Texte hervorheben, im Browser sind die wichtigen W\u00f6rter fett angestrichen.\",\n    \"4. Mit Style Attribut kann man die wichtigen W\u00f6rter markieren, im Browser sind die wichtigen W\u00f6rter fett angestrichen.\"\n  ],\n  [\n    \"5.Mit 1x strong kann ich das Wichtigste markieren, im Browser sind die markierten W\u00f6rter fett angestrichen.\",\n    \"5.Mit HTML Format Element kann ich man das die wichtigen  W\u00f6rter markieren, im Browser sind die markierten W\u00f6rter fett angestrichen.\"\n  ]\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  for (const [oldText, newText] of replacements) {\n    if (para.text === oldText) {\n      para.insertText(newText, Word.InsertLocation.replace);\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# Applies the text changes described by the diff:\n#  - \"Word Wide Web General\" bullets 1-3 reworded\n#  - \"Web Server\" bullet 1 extended, bullet 4 reworded\n#  - \"Web Client\" bullet 3 extended\n#  - \"Hello World\" bullets 4 & 5 reworded\n#\n# Find/Replace (rather than Paragraph.Range.Text =) is used because a\n# straight Range.Text assignment only stomps the first run of a\n# multi-run paragraph; Find.Execute correctly replaces the whole match\n# regardless of how many runs it spans.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"1.HTTP,CSS und Javascript haben WWW zum Erfolg gemacht.\"; New = \"1.HTTP, HTML und URL haben WWW zum Erfolg gemacht.\" },\n  @{ Old = \"2.IP-Adressen,URI oder DNS adressieren Ressourcen im WEB .\"; New = \"2.URL adressieren Ressourcen im WEB .\" },\n  @{ Old = \"3.Bestandteile einer URL: Verzeichnis,Dateiname,Textmarke,Username,Portnummer,CGI-Nummer,Protkoll und Servername.\"; New = \"3.Bestandteile einer URL: schema, special part, username, host, port, part, query, fragment\" },\n  @{ Old = \"1.Als Webserver bezeichnet man jene Server, die zur Verbreitung von Webinhalten im Internet dienen.\"; New = \"1.Als Webserver bezeichnet man jene Server, die Daten speichern und die zur Verbreitung von Webinhalten im Internet dienen.\" },\n  @{ Old = \"4.Durch ASP,PHP werden Inhalte in Webseiten beschrieben.\"; New = \"4.Durch HTML, PHP, ASHP werden Inhalte in Webseiten beschrieben.\" },\n  @{ Old = \"3.Web Client und Webbrowser sind das gleiche.\"; New = \"3.Web Client und Webbrowser sind nicht das gleiche.\" },\n  @{ Old = \"4 Mit 2x strong kann man Texte hervorheben, im Browser sind die wichtigen W\u00f6rter fett angestrichen.\"; New = \"4. Mit Style Attribut kann man die wichtigen W\u00f6rter markieren, im Browser sind die wichtigen W\u00f6rter fett angestrichen.\" },\n  @{ Old = \"5.Mit 1x strong kann ich das Wichtigste markieren, im Browser sind die markierten W\u00f6rter fett angestrichen.\"; New = \"5.Mit HTML Format Element kann ich man das die wichtigen  W\u00f6rter markieren, im Browser sind die markierten W\u00f6rter fett angestrichen.\" }\n)\n\nforeach ($rep in $replacements) {\n  $rng = $d.Content\n  $rng.Find.ClearFormatting()\n  $rng.Find.Replacement.ClearFormatting()\n  $rng.Find.Execute($rep.Old, $false, $true, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null\n}\n"}
